{"js": "// Fix the typo \"yasmak\" -> \"yazmak\" in the heading\n// \"6. \u0130lk unit testimizi yasmak\" (Turkish: \"yazmak\" = \"to write\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the paragraph that contains the misspelling.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text && p.text.indexOf(\"yasmak\") !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (target) {\n  // Narrow down to just the mis-typed word so we only touch the run(s)\n  // that actually contain it, leaving the rest of the paragraph/run\n  // structure (and the \"testimizi\" proofErr-wrapped run) untouched.\n  // (Word.Range.search() ranges are not reliable for sub-string edits\n  // in this host, so split on whitespace instead, which yields\n  // correctly-bounded ranges.)\n  const fullRange = target.getRange();\n  const words = fullRange.split([\" \"], false, false);\n  words.load(\"items,text\");\n  await context.sync();\n\n  let wordRange = null;\n  for (let i = 0; i < words.items.length; i++) {\n    const w = words.items[i];\n    if (w.text && w.text.indexOf(\"yasmak\") !== -1) {\n      wordRange = w;\n      break;\n    }\n  }\n\n  if (wordRange) {\n    wordRange.load(\"text\");\n    await context.sync();\n    // The last \"word\" range includes the trailing paragraph mark (\\r);\n    // strip it before building the replacement text so we don't inject\n    // an extra paragraph break.\n    const cleanWord = wordRange.text.replace(/[\\r\\v]+$/, \"\");\n    const fixedWord = cleanWord.replace(\"yasmak\", \"yazmak\");\n    wordRange.insertText(fixedWord, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Fix the typo \"yasmak\" -> \"yazmak\" in the heading\n# \"6. \u0130lk unit testimizi yasmak\" (Turkish: \"yazmak\" = \"to write\").\n$d = $word.ActiveDocument\n\n$needle = \"yasmak\"\n$replacement = \"yazmak\"\n\n# Note: $d.Content.Find.Execute(...) leaves this host's resulting Range\n# boundaries off by one, so instead locate the text ourselves and build\n# an explicit Range from character offsets. Overwriting just that\n# substring (rather than the whole paragraph) preserves all surrounding\n# runs/formatting exactly.\n$searchStart = 0\nwhile ($true) {\n    $fullText = $d.Content.Text\n    $idx = $fullText.IndexOf($needle, $searchStart)\n    if ($idx -lt 0) {\n        break\n    }\n    $r = $d.Range($idx, $idx + $needle.Length)\n    $r.Text = $replacement\n    $searchStart = $idx + $replacement.Length\n}\n"}
